$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.779.65"
$ws.Range("E2").Value = "  +2.10%  "

$ws.Range("D3").Value = "1.874.73"
$ws.Range("E3").Value = "  +2.26%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'326.77"
$ws.Range("E5").Value = "  -0.56%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").Value = "'0.4652"
$ws.Range("E7").Value = "  +0.71%  "

$ws.Range("D8").Value = "'0.3939"
$ws.Range("E8").Value = "  +2.37%  "

$ws.Range("D9").Value = "'0.07905"
$ws.Range("E9").Value = "  +0.92%  "

$ws.Range("D10").Value = "'0.9732"
$ws.Range("E10").Value = "  +1.88%  "

$ws.Range("D11").Value = "'22.33"
$ws.Range("E11").Value = "  +2.18%  "

$ws.Range("D12").Value = "1.899.76"
$ws.Range("E12").Value = "  +4.10%  "

$ws.Range("D13").Value = "'5.751"
$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("D14").Value = "'6.951"
$ws.Range("E14").Value = "  +1.04%  "

$ws.Range("D15").Value = "'0.06994"
$ws.Range("E15").Value = "  +2.07%  "

$ws.Range("D16").Value = "'88.55"
$ws.Range("E16").Value = "  +2.27%  "

$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").Value = "'16.96"

$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("D21").Value = "28.815.38"
$ws.Range("E21").Value = "  +2.11%  "

$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").Value = "'11.10"
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").Value = "2.058.42"
$ws.Range("E25").Value = "  +0.51%  "

$ws.Range("D26").Value = "'153.60"
$ws.Range("E26").Value = "  +0.44%  "

$ws.Range("D27").Value = "'19.40"
$ws.Range("E27").Value = "  +1.16%  "

$ws.Range("D28").Value = "'5.765"
$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("D29").Value = "'2.008"
$ws.Range("E29").Value = "  +1.90%  "

$ws.Range("D30").Value = "'119.63"
$ws.Range("E30").Value = "  +2.43%  "

$ws.Range("D31").Value = "'0.09378"

$ws.Range("D32").Value = "'0.9416"

$ws.Range("D33").Value = "'5.319"
$ws.Range("E33").Value = "  +1.16%  "

$ws.Range("D34").Value = "'1.349"
$ws.Range("E34").Value = "  +2.83%  "

$ws.Range("D35").Value = "'3.348"
$ws.Range("E35").Value = "  -2.77%  "

$ws.Range("E36").Value = "  -1.91%  "

$ws.Range("D37").Value = "'0.02115"
$ws.Range("E37").Value = "  -1.57%  "

$ws.Range("D38").Value = "'1.148"
$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").Value = "'7.939"
$ws.Range("E39").Value = "  +4.66%  "

$ws.Range("D40").Value = "'0.5689"
$ws.Range("E40").Value = "  +1.62%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.1788"
$ws.Range("E41").Value = "  +0.99%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'9.949"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("D43").Value = "'0.07236"
$ws.Range("E43").Value = "  +3.21%  "

$ws.Range("D44").Value = "'11.78"
$ws.Range("E44").Value = "  +1.85%  "

$ws.Range("D45").Value = "'0.5326"
$ws.Range("E45").Value = "  +1.20%  "

$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.134"
$ws.Range("E46").Value = "  -9.34%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'1.853"
$ws.Range("E47").Value = "  +1.30%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'2.109"
$ws.Range("E48").Value = "  -5.59%  "

$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("D50").Value = "'2.364"
$ws.Range("E50").Value = "  +1.65%  "

$ws.Range("E51").Value = "  +0.27%  "
